$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 122 (pushes existing rows 122:133 down to 123:134),
# matching the weekly data refresh described in the commit message.
$ws.Rows("122:122").Insert()

$ws.Range("A122").Value = 4
$ws.Range("B122").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C122").Value = "Los Lagos"
$ws.Range("D122").Value = 44504
$ws.Range("E122").Value = 10
$ws.Range("F122").Value = 100112028
$ws.Range("G122").Value = "Sandia"
$ws.Range("H122").Value = "Sin especificar"
$ws.Range("I122").Value = "Primera"
$ws.Range("J122").Value = 400
$ws.Range("K122").Value = 1000
$ws.Range("L122").Value = 1200
$ws.Range("M122").Value = 1100
$ws.Range("N122").Value = "$/kilo (volumen en unidades)"
$ws.Range("O122").Value = "Perú"
$ws.Range("P122").Value = 1100
$ws.Range("Q122").Value = 1
$ws.Range("R122").Value = "Hortaliza"
